$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Row 76
$ws.Range("H76").Value = 1649.5
$ws.Range("I76").Value = 1199.6666
$ws.Range("K76").Value = 1199.6666
$ws.Range("M76").Value = -884.6666
# Row 79
$ws.Range("H79").Value = 1649.5
$ws.Range("I79").Value = 1199.6666
$ws.Range("K79").Value = 1199.6666
$ws.Range("M79").Value = -107.6666
# Row 116
$ws.Range("H116").Value = 3891.5386
$ws.Range("I116").Value = 3873.818
$ws.Range("K116").Value = 3873.818
$ws.Range("M116").Value = -431.8180000000002
# Row 123
$ws.Range("H123").Value = 164500
$ws.Range("J123").Value = 164500
$ws.Range("L123").Value = 164500
$ws.Range("N123").Value = -174300
# Row 137
$ws.Range("H137").Value = 2738.8
$ws.Range("J137").Value = 4173.25
$ws.Range("L137").Value = 12519.75
$ws.Range("N137").Value = -17619.75
# Row 138
$ws.Range("H138").Value = 6790.108
$ws.Range("J138").Value = 7508.9033
$ws.Range("L138").Value = 22526.7099
$ws.Range("N138").Value = -32806.7099

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 2992.3542
$ws.Range("I32").Value = 2356.3333
$ws.Range("J32").Value = 12532.667
$ws.Range("K32").Value = 2356.3333
$ws.Range("L32").Value = 12532.667
$ws.Range("M32").Value = -2069.3333
$ws.Range("N32").Value = -13106.667
# Row 61
$ws.Range("H61").Value = 3391.8333
$ws.Range("I61").Value = 3114.875
$ws.Range("J61").Value = 3945.75
$ws.Range("K61").Value = 3114.875
$ws.Range("L61").Value = 3945.75
$ws.Range("M61").Value = -2902.875
$ws.Range("N61").Value = -4369.75
# Row 110
$ws.Range("H110").Value = 2750.25
$ws.Range("J110").Value = 909.5
$ws.Range("L110").Value = 909.5
$ws.Range("N110").Value = -4999.5
# Row 132
$ws.Range("H132").Value = 2862.7778
$ws.Range("I132").Value = 2404.1
$ws.Range("J132").Value = 3436.125
$ws.Range("K132").Value = 7212.299999999999
$ws.Range("L132").Value = 10308.375
$ws.Range("M132").Value = -4682.299999999999
$ws.Range("N132").Value = -15368.375
# Row 136
$ws.Range("H136").Value = 3391.8333
$ws.Range("I136").Value = 3114.875
$ws.Range("J136").Value = 3945.75
$ws.Range("K136").Value = 9344.625
$ws.Range("L136").Value = 11837.25
$ws.Range("M136").Value = -6794.625
$ws.Range("N136").Value = -16937.25

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
# Row 134
$ws.Range("H134").Value = 4721.5
$ws.Range("I134").Value = 4721.5
$ws.Range("K134").Value = 14164.5
$ws.Range("M134").Value = -11629.5

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
# Row 22
$ws.Range("H22").Value = 2049
$ws.Range("J22").Value = 0
$ws.Range("L22").Value = 0
$ws.Range("N22").ClearContents()
# Row 31
$ws.Range("H31").Value = 5222.222
$ws.Range("J31").Value = 6783.067
$ws.Range("L31").Value = 6783.067
$ws.Range("N31").Value = -7373.067
# Row 34
$ws.Range("H34").Value = 5222.222
$ws.Range("J34").Value = 6783.067
$ws.Range("L34").Value = 6783.067
$ws.Range("N34").Value = -7187.067
# Row 58
$ws.Range("H58").Value = 2374
$ws.Range("I58").Value = 2365
$ws.Range("K58").Value = 2365
$ws.Range("M58").Value = -2162
# Row 62
$ws.Range("H62").Value = 0
$ws.Range("J62").Value = 0
$ws.Range("L62").Value = 0
$ws.Range("N62").ClearContents()
# Row 65
$ws.Range("H65").Value = 0
$ws.Range("J65").Value = 0
$ws.Range("L65").Value = 0
$ws.Range("N65").ClearContents()
# Row 132
$ws.Range("H132").Value = 4511
$ws.Range("I132").Value = 3972
$ws.Range("J132").Value = 5499.1665
$ws.Range("K132").Value = 11916
$ws.Range("L132").Value = 16497.4995
$ws.Range("M132").Value = -9386
$ws.Range("N132").Value = -21557.4995
# Row 136
$ws.Range("H136").Value = 2374
$ws.Range("I136").Value = 2365
$ws.Range("K136").Value = 7095
$ws.Range("M136").Value = -4545

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
# Row 92
$ws.Range("H92").Value = 3
$ws.Range("J92").Value = 3
$ws.Range("L92").Value = 9
$ws.Range("N92").Value = -2505
# Row 107
$ws.Range("H107").Value = 2931.6667
$ws.Range("J107").Value = 954.6667
$ws.Range("L107").Value = 2864.0001
$ws.Range("N107").Value = -6704.0001
# Row 128
$ws.Range("H128").Value = 249999.5
$ws.Range("I128").Value = 249999.5
$ws.Range("K128").Value = 749998.5
$ws.Range("M128").Value = -745018.5
# Row 131
$ws.Range("H131").Value = 1300.4445
$ws.Range("I131").Value = 602.75
$ws.Range("J131").Value = 1499.7858
$ws.Range("K131").Value = 1808.25
$ws.Range("L131").Value = 4499.357400000001
$ws.Range("M131").Value = 3231.75
$ws.Range("N131").Value = -14579.3574

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
# Row 102
$ws.Range("H102").Value = 1533.9
$ws.Range("I102").Value = 1602.7778
$ws.Range("K102").Value = 1602.7778
$ws.Range("M102").Value = 19.22219999999993
# Row 119
$ws.Range("H119").Value = 65000
$ws.Range("J119").Value = 65000
$ws.Range("L119").Value = 65000
$ws.Range("N119").Value = -74676
# Row 132
$ws.Range("H132").Value = 5202.857
$ws.Range("I132").Value = 4985.8184
$ws.Range("J132").Value = 5998.6665
$ws.Range("K132").Value = 14957.4552
$ws.Range("L132").Value = 17995.9995
$ws.Range("M132").Value = -12427.4552
$ws.Range("N132").Value = -23055.9995

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Range("H7").Value = 3641.7144
$ws.Range("I7").Value = 2832.3333
$ws.Range("J7").Value = 4248.75
$ws.Range("K7").Value = 2832.3333
$ws.Range("L7").Value = 4248.75
$ws.Range("M7").Value = -2720.3333
$ws.Range("N7").Value = -4472.75
# Row 10
$ws.Range("H10").Value = 2401.4
$ws.Range("I10").Value = 2003
$ws.Range("K10").Value = 2003
$ws.Range("M10").Value = -1863
# Row 126
$ws.Range("H126").Value = 3641.7144
$ws.Range("I126").Value = 2832.3333
$ws.Range("J126").Value = 4248.75
$ws.Range("K126").Value = 8496.999899999999
$ws.Range("L126").Value = 12746.25
$ws.Range("M126").Value = -6026.999899999999
$ws.Range("N126").Value = -17686.25
# Row 132
$ws.Range("H132").Value = 5748
$ws.Range("I132").Value = 5004
$ws.Range("J132").Value = 5996
$ws.Range("K132").Value = 15012
$ws.Range("L132").Value = 17988
$ws.Range("M132").Value = -12482
$ws.Range("N132").Value = -23048
# Row 136
$ws.Range("H136").Value = 25525.428
$ws.Range("I136").Value = 3005.5
$ws.Range("K136").Value = 9016.5
$ws.Range("M136").Value = -6466.5

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
# Row 122
$ws.Range("H122").Value = 499
$ws.Range("I122").Value = 499
$ws.Range("K122").Value = 1497
$ws.Range("M122").Value = 953
# Row 126
$ws.Range("H126").Value = 1942.9
$ws.Range("I126").Value = 1812
$ws.Range("J126").Value = 2248.3333
$ws.Range("K126").Value = 5436
$ws.Range("L126").Value = 6744.999899999999
$ws.Range("M126").Value = -2966
$ws.Range("N126").Value = -11684.9999
# Row 132
$ws.Range("H132").Value = 2811.5
$ws.Range("I132").Value = 2509.6667
$ws.Range("K132").Value = 7529.000100000001
$ws.Range("M132").Value = -4999.000100000001
# Row 136
$ws.Range("H136").Value = 11765.667
$ws.Range("I136").Value = 14478.223
$ws.Range("J136").Value = 3628
$ws.Range("K136").Value = 43434.669
$ws.Range("L136").Value = 10884
$ws.Range("M136").Value = -40884.669
$ws.Range("N136").Value = -15984
